# Reorder test files: OBP, SLG
#
# The workbook has four "game" sheets (1-4), each with stat columns laid
# out as ... N=AVG, O=SLG, P=OBP, Q=OPS. This swaps O and P so the order
# becomes N=AVG, O=OBP, P=SLG, Q=OPS (OBP before SLG).

$wb = $excel.ActiveWorkbook

# Workbook-level view: the saved file no longer pins a specific
# first-visible-tab / active-tab pair on the workbookView element.
$wb.Windows.Item(1).DisplayWorkbookTabs = $true

$sheetNames = @("Fall 2015 09.16", "Fall 2015 09.09", "Spring 2014 04.16", "Spring 2014 04.09")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # --- swap the OBP/SLG columns (O <-> P) for header + rows 2-6 ---
    for ($row = 2; $row -le 6; $row++) {
        $oFormula = $ws.Range("O$row").Formula
        $pFormula = $ws.Range("P$row").Formula
        $ws.Range("O$row").Formula = $pFormula
        $ws.Range("P$row").Formula = $oFormula
    }

    $ws.Range("O1").Value = "OBP"
    $ws.Range("P1").Value = "SLG"

    # Q column now reads P2+O2 instead of O2+P2 (same result, new formula text)
    for ($row = 2; $row -le 6; $row++) {
        $ws.Range("Q$row").Formula = "=IF(P$row,P$row+O$row,)"
    }

    # row 6 totals (B6:L6) re-entered as one fill so they share a single formula
    $ws.Range("B6").Formula = "=SUM(B2:B5)"
    $ws.Range("C6").Formula = "=SUM(C2:C5)"
    $ws.Range("D6").Formula = "=SUM(D2:D5)"
    $ws.Range("E6").Formula = "=SUM(E2:E5)"
    $ws.Range("F6").Formula = "=SUM(F2:F5)"
    $ws.Range("G6").Formula = "=SUM(G2:G5)"
    $ws.Range("H6").Formula = "=SUM(H2:H5)"
    $ws.Range("I6").Formula = "=SUM(I2:I5)"
    $ws.Range("J6").Formula = "=SUM(J2:J5)"
    $ws.Range("K6").Formula = "=SUM(K2:K5)"
    $ws.Range("L6").Formula = "=SUM(L2:L5)"

    # selection moved to the (now-OBP) column O, whole-column selected
    $ws.Range("O1:O1048576").Select()

    $ws.PageSetup.Orientation = 1
}

# Tab selection moves from "Spring 2014 04.09" (4th sheet) to
# "Fall 2015 09.16" (1st sheet).
$wb.Worksheets.Item("Fall 2015 09.16").Activate()
